$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose "Beteckning"/"Datum"/"Area (ha)" data got reshuffled (rows 7-14) ---
# New row 7 <- old row 9
$ws.Range("A7").Value = "A 12651-2022"
$ws.Range("B7").Value = 44641
$ws.Range("G7").Value = 3.2

# New row 8 <- old row 12
$ws.Range("A8").Value = "A 8194-2025"
$ws.Range("B8").Value = 45708
$ws.Range("G8").Value = 1.9

# New row 9 <- old row 10
$ws.Range("A9").Value = "A 5792-2024"
$ws.Range("B9").Value = 45335
$ws.Range("G9").Value = 5.6

# New row 10 <- old row 14
$ws.Range("A10").Value = "A 50997-2025"
$ws.Range("B10").Value = 45946
$ws.Range("G10").Value = 1.5

# New row 11 <- old row 7
$ws.Range("A11").Value = "A 7827-2026"
$ws.Range("B11").Value = 46062.63958333333
$ws.Range("G11").Value = 2.1

# New row 12 <- old row 8
$ws.Range("A12").Value = "A 7814-2026"
$ws.Range("B12").Value = 46062.61388888889
$ws.Range("G12").Value = 1.1

# New row 13 <- old row 11
$ws.Range("A13").Value = "A 13651-2023"
$ws.Range("B13").Value = 45006
$ws.Range("G13").Value = 2.2

# New row 14 <- old row 13
$ws.Range("A14").Value = "A 35642-2023"
$ws.Range("B14").Value = 45147
$ws.Range("G14").Value = 1.2

# --- "Förändrad" (column C) refresh date: 46066 -> 46070 for every data row (2-16) ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Range("C$r").Value = 46070
}
